# ClinicalData_template.xlsx edit
# Rename the "tissue id" / "intervention id" / "disease id" header columns
# (B1:D1) to "tissue" / "intervention" / "disease" — these columns are now
# populated with the CKG internal identifiers instead of the raw ids, so
# the " id" suffix is dropped from the column header label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "tissue"
$ws.Range("C1").Value = "intervention"
$ws.Range("D1").Value = "disease"

# Re-apply the "best fit" column widths that Excel recalculated for every
# data column (B:AO) after the header text changed. Excel's ColumnWidth
# property is expressed in characters; the stored OOXML width include a
# fixed ~0.83 (5/6) character padding that Excel adds automatically, so we
# back that out here so the persisted width matches the target exactly.
$ws.Columns.Item(2).ColumnWidth = 6.666666666666667
$ws.Columns.Item(3).ColumnWidth = 12.0
$ws.Columns.Item(4).ColumnWidth = 7.833333333333333
$ws.Columns.Item(5).ColumnWidth = 23.666666666666668
$ws.Columns.Item(6).ColumnWidth = 21.333333333333332
$ws.Columns.Item(7).ColumnWidth = 26.166666666666668
$ws.Columns.Item(8).ColumnWidth = 19.833333333333332
$ws.Columns.Item(9).ColumnWidth = 27.166666666666668
$ws.Columns.Item(10).ColumnWidth = 34.166666666666664
$ws.Columns.Item(11).ColumnWidth = 20.333333333333332
$ws.Columns.Item(12).ColumnWidth = 19.166666666666668
$ws.Columns.Item(13).ColumnWidth = 23.5
$ws.Columns.Item(14).ColumnWidth = 21.166666666666668
$ws.Columns.Item(15).ColumnWidth = 26.0
$ws.Columns.Item(16).ColumnWidth = 19.666666666666668
$ws.Columns.Item(17).ColumnWidth = 27.0
$ws.Columns.Item(18).ColumnWidth = 34.0
$ws.Columns.Item(19).ColumnWidth = 20.166666666666668
$ws.Columns.Item(20).ColumnWidth = 19.0
$ws.Columns.Item(21).ColumnWidth = 8.166666666666666
$ws.Columns.Item(22).ColumnWidth = 12.333333333333334
$ws.Columns.Item(23).ColumnWidth = 8.166666666666666
$ws.Columns.Item(24).ColumnWidth = 8.166666666666666
$ws.Columns.Item(25).ColumnWidth = 9.833333333333334
$ws.Columns.Item(26).ColumnWidth = 12.833333333333334
$ws.Columns.Item(27).ColumnWidth = 16.5
$ws.Columns.Item(28).ColumnWidth = 13.5
$ws.Columns.Item(29).ColumnWidth = 12.0
$ws.Columns.Item(30).ColumnWidth = 27.5
$ws.Columns.Item(31).ColumnWidth = 22.0
$ws.Columns.Item(32).ColumnWidth = 30.0
$ws.Columns.Item(33).ColumnWidth = 9.333333333333334
$ws.Columns.Item(34).ColumnWidth = 19.0
$ws.Columns.Item(35).ColumnWidth = 43.166666666666664
$ws.Columns.Item(36).ColumnWidth = 40.666666666666664
$ws.Columns.Item(37).ColumnWidth = 5.5
$ws.Columns.Item(38).ColumnWidth = 5.833333333333333
$ws.Columns.Item(39).ColumnWidth = 6.5
$ws.Columns.Item(40).ColumnWidth = 7.666666666666667
$ws.Columns.Item(41).ColumnWidth = 8.5

# Move the view so it scrolls to the right-hand side of the new, wider
# sheet (to column R) and leave the active selection on E1.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 18
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E1").Select()
